$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(117, 2).Value = 7013885
$ws.Cells.Item(117, 29).Value2 = -0
$ws.Cells.Item(118, 27).Value2 = 0.8999999999999999
Write-Output "done"
